$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect it so the cells below can be updated.
$ws.Unprotect()

# Update the confidential note text (date change 2021-05-07 -> 2021-05-10)
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-10 for illustrative purposes only and are subject to change."
# Undo the automatic row-height bump caused by the embedded line break so the
# row keeps using the sheet's default (un-customized) height.
$ws.Rows("10").AutoFit()

# Update the Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.2460498812741683
$ws.Range("E2").Value = -0.01977261492832416

$ws.Range("D3").Value = 0.4990156225939272
$ws.Range("E3").Value = 0.001558846453624252

$ws.Range("D4").Value = 0.09561167178831745
$ws.Range("E4").Value = -0.01608789484010198

$ws.Range("D5").Value = 0.1024603548848817
$ws.Range("E5").Value = -0.0006533333333333946

$ws.Range("D6").Value = 0.05686246945870525
$ws.Range("E6").Value = -0.01934623082054698

$ws.Range("D7").Value = 0.9999999999999999
$ws.Range("E7").Value = -0.006792366567593078

# Restore protection on the sheet.
$ws.Protect()
